$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Swap-RowRange {
    param(
        [object]$Worksheet,
        [int]$Row1,
        [int]$Row2,
        [string]$StartCol,
        [string]$EndCol
    )

    $range1 = $Worksheet.Range("$StartCol$Row1`:$EndCol$Row1")
    $range2 = $Worksheet.Range("$StartCol$Row2`:$EndCol$Row2")

    $values1 = $range1.Value2
    $values2 = $range2.Value2

    $range1.Value2 = $values2
    $range2.Value2 = $values1
}

# Rows 20 and 21 swap their B:AC content (row id in column A stays put).
Swap-RowRange $ws 20 21 "B" "AC"

# Rows 117 and 118 swap their B:AC content (row id in column A stays put).
Swap-RowRange $ws 117 118 "B" "AC"
